{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the contract-detail changes described by the diff:\n//   1. Worker identity block: name/CC/address swapped to the new worker.\n//   2. Job title (clause PRIMERA): \"Asistente de Recursos Humanos\" -> \"Desarrollador Web\".\n//   3. Contract start date: 2019-05-10 -> 2021-09-01.\n//   4. Job title (clause TERCERA): \"Asistente de Recursos Humanos\" -> \"Desarrollador Web\".\n//   5. Monthly salary: $1800000 -> $3200000.\n//   6. Worker's first-name-only signature line: \"Mar\u00eda\" -> \"Andr\u00e9s\".\n//\n// Each replacement is done via Body.search() + Range.insertText(text, \"Replace\"),\n// using long, unique search strings so only the intended run is touched.\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, newText, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${searchText}`);\n  }\n\n  // Replace every match (the diff's two \"Asistente de Recursos Humanos\"\n  // occurrences are handled by two separate, distinctly-worded calls below,\n  // so in practice each call here hits exactly one paragraph).\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n// 1) Worker identity line (name, CC number, address) in the intro paragraph.\n//    NOTE: the document uses plain straight quotes (\") around \"El Trabajador\", not curly quotes.\nawait replaceOnce(\n  \"Rodr\u00edguez L\u00f3pez Mar\u00eda, identificado/a con CC 87654321, domiciliado/a en Carrera 7 #32-18, Medell\u00edn, quien en adelante se denominar\u00e1 \\\"El Trabajador\\\",\",\n  \"Gonz\u00e1lez Mora Andr\u00e9s, identificado/a con CC 11223344, domiciliado/a en Av. Las Am\u00e9ricas #45-67, Cali, quien en adelante se denominar\u00e1 \\\"El Trabajador\\\",\"\n);\n\n// 2) Job title mentioned in clause PRIMERA (natural of the contract).\nawait replaceOnce(\n  \"El Empleador contrata los servicios de El Trabajador para desempe\u00f1ar el cargo de Asistente de Recursos Humanos, bajo las condiciones estipuladas en este contrato y las disposiciones legales aplicables.\",\n  \"El Empleador contrata los servicios de El Trabajador para desempe\u00f1ar el cargo de Desarrollador Web, bajo las condiciones estipuladas en este contrato y las disposiciones legales aplicables.\"\n);\n\n// 3) Contract start date in clause SEGUNDA.\nawait replaceOnce(\n  \"Tipo de contrato Indefinido , comenzando el d\u00eda 2019-05-10\",\n  \"Tipo de contrato Indefinido , comenzando el d\u00eda 2021-09-01\"\n);\n\n// 4) Job title mentioned in clause TERCERA (place/schedule of work).\nawait replaceOnce(\n  \"El Trabajador desempe\u00f1ar\u00e1 sus funciones en Asistente de Recursos Humanos .\",\n  \"El Trabajador desempe\u00f1ar\u00e1 sus funciones en Desarrollador Web .\"\n);\n\n// 5) Monthly salary amount in clause CUARTA.\nawait replaceOnce(\n  \"El Trabajador recibir\u00e1 un salario mensual de $1800000, que se pagar\u00e1 de manera mensual, sujeto a las deducciones legales correspondientes.\",\n  \"El Trabajador recibir\u00e1 un salario mensual de $3200000, que se pagar\u00e1 de manera mensual, sujeto a las deducciones legales correspondientes.\"\n);\n\n// 6) Worker's given name on the signature line (standalone \"Mar\u00eda\" run).\nawait replaceOnce(\"Mar\u00eda\", \"Andr\u00e9s\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the contract-detail changes described by the diff:\n#   1. Worker identity block: name/CC/address swapped to the new worker.\n#   2. Job title (clause PRIMERA): \"Asistente de Recursos Humanos\" -> \"Desarrollador Web\".\n#   3. Contract start date: 2019-05-10 -> 2021-09-01.\n#   4. Job title (clause TERCERA): \"Asistente de Recursos Humanos\" -> \"Desarrollador Web\".\n#   5. Monthly salary: $1800000 -> $3200000.\n#   6. Worker's first-name-only signature line: \"Mar\u00eda\" -> \"Andr\u00e9s\".\n#\n# Each change is located with Range.Find.Execute() (search only, no Replace\n# argument) and then applied by assigning Range.Text directly. Doing the\n# substitution this way \u2014 instead of passing a Replace:= string to\n# Find.Execute \u2014 avoids Word's Find/Replace \"smart quotes\" autocorrection,\n# which would otherwise turn the straight double quotes around\n# \"El Trabajador\" into curly ones.\n\n$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch([string]$searchText, [string]$newText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.MatchCase = $true\n    $find.MatchWildcards = $false\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Search text not found: $searchText\"\n    }\n    $rng.Text = $newText\n}\n\n# 1) Worker identity line (name, CC number, address) in the intro paragraph.\nReplace-FirstMatch `\n    \"Rodr\u00edguez L\u00f3pez Mar\u00eda, identificado/a con CC 87654321, domiciliado/a en Carrera 7 #32-18, Medell\u00edn, quien en adelante se denominar\u00e1 `\"El Trabajador`\",\" `\n    \"Gonz\u00e1lez Mora Andr\u00e9s, identificado/a con CC 11223344, domiciliado/a en Av. Las Am\u00e9ricas #45-67, Cali, quien en adelante se denominar\u00e1 `\"El Trabajador`\",\"\n\n# 2) Job title mentioned in clause PRIMERA (natural of the contract).\nReplace-FirstMatch `\n    \"El Empleador contrata los servicios de El Trabajador para desempe\u00f1ar el cargo de Asistente de Recursos Humanos, bajo las condiciones estipuladas en este contrato y las disposiciones legales aplicables.\" `\n    \"El Empleador contrata los servicios de El Trabajador para desempe\u00f1ar el cargo de Desarrollador Web, bajo las condiciones estipuladas en este contrato y las disposiciones legales aplicables.\"\n\n# 3) Contract start date in clause SEGUNDA.\nReplace-FirstMatch `\n    \"Tipo de contrato Indefinido , comenzando el d\u00eda 2019-05-10\" `\n    \"Tipo de contrato Indefinido , comenzando el d\u00eda 2021-09-01\"\n\n# 4) Job title mentioned in clause TERCERA (place/schedule of work).\nReplace-FirstMatch `\n    \"El Trabajador desempe\u00f1ar\u00e1 sus funciones en Asistente de Recursos Humanos .\" `\n    \"El Trabajador desempe\u00f1ar\u00e1 sus funciones en Desarrollador Web .\"\n\n# 5) Monthly salary amount in clause CUARTA.\nReplace-FirstMatch `\n    \"El Trabajador recibir\u00e1 un salario mensual de `$1800000, que se pagar\u00e1 de manera mensual, sujeto a las deducciones legales correspondientes.\" `\n    \"El Trabajador recibir\u00e1 un salario mensual de `$3200000, que se pagar\u00e1 de manera mensual, sujeto a las deducciones legales correspondientes.\"\n\n# 6) Worker's given name on the signature line (standalone \"Mar\u00eda\" run).\n#    By this point the only remaining \"Mar\u00eda\" is the signature-block one,\n#    since step 1 already rewrote the other occurrence.\nReplace-FirstMatch \"Mar\u00eda\" \"Andr\u00e9s\"\n"}
